# FaQSheet.xlsx edit: append 10 new FAQ rows (rows 64-73) to sheet "1",
# matching the shared-strings / sheetData additions from the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting: reuse existing border-only styles by copy/paste-special -------
# Rows 64-69, column A carries the same "Client Requisition" grouping style that
# row 63's A/B cells already use; columns B/C (and all of rows 70-73) use the
# plain bordered style already used throughout the table (e.g. B61).
for ($r = 64; $r -le 69; $r++) {
  $ws.Cells.Item(63,1).Copy()
  $ws.Cells.Item($r,1).PasteSpecial(-4122)
  $ws.Cells.Item(61,2).Copy()
  $ws.Cells.Item($r,2).PasteSpecial(-4122)
  $ws.Cells.Item($r,3).PasteSpecial(-4122)
}
for ($r = 70; $r -le 73; $r++) {
  $ws.Cells.Item(61,2).Copy()
  $ws.Cells.Item($r,1).PasteSpecial(-4122)
  $ws.Cells.Item($r,2).PasteSpecial(-4122)
  $ws.Cells.Item($r,3).PasteSpecial(-4122)
}

# --- Row 64: What mode will a construction project be executed on ? -----------
$ws.Cells.Item(64,1).Value = "Client Requisition"
$ws.Cells.Item(64,2).Value = "What mode will a construction project be executed on ?"
$ws.Cells.Item(64,3).Value = "Geographical mode (<100 crores) or Project mode (>100 crores or important client or project)"

# --- Row 65: team composition exception ----------------------------------------
$ws.Cells.Item(65,1).Value = "Client Requisition"
$ws.Cells.Item(65,2).Value = "Is there any exception to team composition in a Project Mode assignment ?"
$ws.Cells.Item(65,3).Value = "Decision to form a separate team is also taken by Director General (DG) CPWD considering the parameters like the work load, quantum of work, availability staff etc."

# --- Row 66: Executing Engineer/Office ------------------------------------------
$ws.Cells.Item(66,1).Value = "Client Requisition"
$ws.Cells.Item(66,2).Value = "Who should be the Executing Engineer/Office for a project?"
$ws.Cells.Item(66,3).Value = "If Geographical mode, concerned Division office/Circle Office/Zonal Office in the geographical jurisdiction`nIf Project Mode, Project Manager (PM) or Chief Project Manager (CPM)"

# --- Row 67: team composition for Project Mode ----------------------------------
$ws.Cells.Item(67,1).Value = "Client Requisition"
$ws.Cells.Item(67,2).Value = "What is team composition for a Project Mode assignment ?"
$ws.Cells.Item(67,3).Value = "PM/CPM, EE, AE, AEE and JE"

# --- Row 68: clients -------------------------------------------------------------
$ws.Cells.Item(68,1).Value = "Client Requisition"
$ws.Cells.Item(68,2).Value = "Who are our clients?"
$ws.Cells.Item(68,3).Value = "The clientele can be categorized as follows :- `n1. Ministry of GoI (leading to Budgeted and Authorised Works)`n2. Autonomous bodies (leading to Deposit Works) `n3. PSU's (leading to Deposit Works) `n4. NGO's (leading to Deposit Works) `n5. Private organization (leading to Deposit Works)"

# --- Row 69: Competent Authority --------------------------------------------------
$ws.Cells.Item(69,1).Value = "Client Requisition"
$ws.Cells.Item(69,2).Value = "What is the basis of deciding Competent Authority for any approval ?"
$ws.Cells.Item(69,3).Value = "Competent Authority is decided as per the Delegation of Financial Powers (DFPR) as mentioned in Table A"

# --- Row 70: Preliminary estimates (Table A - Question/Answer header pair) -------
$ws.Cells.Item(70,1).Value = "Preliminary estimates (PE)"
$ws.Cells.Item(70,2).Value = "Question"
$ws.Cells.Item(70,3).Value = "Answer"

# --- Row 71: Administrative approvals and expenditure sanction -------------------
$ws.Cells.Item(71,1).Value = "Administrative approvals and expenditure sanction"
$ws.Cells.Item(71,2).Value = "Question"
$ws.Cells.Item(71,3).Value = "Answer"

# --- Row 72: Detailed estimates ---------------------------------------------------
$ws.Cells.Item(72,1).Value = "Detailed estimates"
$ws.Cells.Item(72,2).Value = "Question"
$ws.Cells.Item(72,3).Value = "Answer"

# --- Row 73: Technical sanction ---------------------------------------------------
$ws.Cells.Item(73,1).Value = "Technical sanction"
$ws.Cells.Item(73,2).Value = "Question"
$ws.Cells.Item(73,3).Value = "Answer"

# --- View state: select the newly-added block and scroll it into view ------------
$ws.Range("A64:A73").Select()
$excel.ActiveWindow.ScrollRow = 60

Write-Output "FaQSheet: added rows 64-73"
